# Append 2025-04-26 price row (row 56) to each Argent/Solar price sheet,
# repeating the last known (2025-04-25) price for every series.

$wb = $excel.ActiveWorkbook

$newDate = "2025-04-26"

$sheetValues = @{
    "N-Dense"                    = "38"
    "N-Type"                     = "37.78"
    "N-type Wafer"                = "1.15"
    "Cell Topcon 183mm"           = "0.278"
    "Module Topcon 183mm"         = "0.09"
    "Silver Rear_side"            = "5,424"
    "Silver Busbar front-side"    = "8,121"
    "Silver finger front-side"    = "8,171"
    "USD_CNY"                     = "7.3083"
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if ($sheetValues.ContainsKey($name)) {
        $price = $sheetValues[$name]

        $rowIndex = 56

        $dateCell = $ws.Cells.Item($rowIndex, 1)
        $dateCell.Value = "'" + $newDate
        $dateCell.Style = "Normal"

        $priceCell = $ws.Cells.Item($rowIndex, 2)
        $priceCell.Value = "'" + $price
        $priceCell.Style = "Normal"
    }
}
